# Auto-generated Excel COM-interop script
# Applies value updates to the Ravana_Profits crafting-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 681.6
$ws.Range("I80").Value = 729
$ws.Range("J80").Value = 610.5
$ws.Range("K80").Value = 2187
$ws.Range("L80").Value = 1831.5
$ws.Range("M80").Value = -1189
$ws.Range("N80").Value = -3827.5

$ws.Range("H83").Value = 681.6
$ws.Range("I83").Value = 729
$ws.Range("J83").Value = 610.5
$ws.Range("K83").Value = 6561
$ws.Range("L83").Value = 5494.5
$ws.Range("M83").Value = -1569
$ws.Range("N83").Value = -15478.5

$ws.Range("H88").Value = 1498.3334
$ws.Range("I88").Value = 1497
$ws.Range("K88").Value = 1497
$ws.Range("M88").Value = -1091

$ws.Range("H91").Value = 1498.3334
$ws.Range("I91").Value = 1497
$ws.Range("K91").Value = 1497
$ws.Range("M91").Value = -93

$ws.Range("H118").Value = 1049.25
$ws.Range("I118").Value = 465.66666
$ws.Range("K118").Value = 1396.99998
$ws.Range("M118").Value = 260.0000199999999

$ws.Range("H129").Value = 1405.7778
$ws.Range("I129").Value = 1358.6666
$ws.Range("K129").Value = 4075.9998
$ws.Range("M129").Value = 924.0001999999999

$ws.Range("H141").Value = 3136.8572
$ws.Range("J141").Value = 1596
$ws.Range("L141").Value = 4788
$ws.Range("N141").Value = -15148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6553.0586
$ws.Range("I32").Value = 6315.2144
$ws.Range("K32").Value = 6315.2144
$ws.Range("M32").Value = -6028.2144

$ws.Range("H61").Value = 2997
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2997
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 2997
$ws.Range("N61").Value = -3421
$ws.Range("M61").ClearContents()

$ws.Range("H74").Value = 1769.3334
$ws.Range("I74").Value = 2960
$ws.Range("K74").Value = 2960
$ws.Range("M74").Value = -2086

$ws.Range("H76").Value = 40000
$ws.Range("J76").Value = 40000
$ws.Range("L76").Value = 40000
$ws.Range("N76").Value = -40676

$ws.Range("H77").Value = 1769.3334
$ws.Range("I77").Value = 2960
$ws.Range("K77").Value = 14800
$ws.Range("M77").Value = -10432

$ws.Range("H79").Value = 40000
$ws.Range("J79").Value = 40000
$ws.Range("L79").Value = 40000
$ws.Range("N79").Value = -42340

$ws.Range("H110").Value = 1340.3529
$ws.Range("I110").Value = 1437.9231
$ws.Range("J110").Value = 1023.25
$ws.Range("K110").Value = 1437.9231
$ws.Range("L110").Value = 1023.25
$ws.Range("M110").Value = 607.0769
$ws.Range("N110").Value = -5113.25

$ws.Range("H136").Value = 2997
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2997
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 8991
$ws.Range("N136").Value = -14091
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 722.9
$ws.Range("I80").Value = 1151.75
$ws.Range("J80").Value = 437
$ws.Range("K80").Value = 1151.75
$ws.Range("L80").Value = 437
$ws.Range("M80").Value = -153.75
$ws.Range("N80").Value = -2433

$ws.Range("H83").Value = 722.9
$ws.Range("I83").Value = 1151.75
$ws.Range("J83").Value = 437
$ws.Range("K83").Value = 5758.75
$ws.Range("L83").Value = 2185
$ws.Range("M83").Value = -766.75
$ws.Range("N83").Value = -12169

$ws.Range("H134").Value = 2567.5293
$ws.Range("I134").Value = 2352.125
$ws.Range("K134").Value = 7056.375
$ws.Range("M134").Value = -4521.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1750.0625
$ws.Range("I31").Value = 1882.1818
$ws.Range("J31").Value = 1459.4
$ws.Range("K31").Value = 1882.1818
$ws.Range("L31").Value = 1459.4
$ws.Range("M31").Value = -1587.1818
$ws.Range("N31").Value = -2049.4

$ws.Range("H34").Value = 1750.0625
$ws.Range("I34").Value = 1882.1818
$ws.Range("J34").Value = 1459.4
$ws.Range("K34").Value = 1882.1818
$ws.Range("L34").Value = 1459.4
$ws.Range("M34").Value = -1680.1818
$ws.Range("N34").Value = -1863.4

$ws.Range("H105").Value = 3359.6
$ws.Range("I105").Value = 2949.5
$ws.Range("K105").Value = 2949.5
$ws.Range("M105").Value = -1202.5

$ws.Range("H132").Value = 3565.8572
$ws.Range("I132").Value = 3393.84
$ws.Range("K132").Value = 10181.52
$ws.Range("M132").Value = -7651.52

$ws.Range("H134").Value = 4350.8335
$ws.Range("I134").Value = 4394.75
$ws.Range("K134").Value = 13184.25
$ws.Range("M134").Value = -10649.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 2748.5
$ws.Range("J23").Value = 3564.6667
$ws.Range("L23").Value = 10694.0001
$ws.Range("N23").Value = -11164.0001

$ws.Range("H34").Value = 2766.8333
$ws.Range("I34").Value = 532.6667
$ws.Range("J34").Value = 5001
$ws.Range("K34").Value = 1598.0001
$ws.Range("L34").Value = 15003
$ws.Range("M34").Value = -1514.0001
$ws.Range("N34").Value = -15171

$ws.Range("H55").Value = 5125
$ws.Range("I55").Value = 1500
$ws.Range("J55").Value = 6333.3335
$ws.Range("K55").Value = 4500
$ws.Range("L55").Value = 19000.0005
$ws.Range("M55").Value = -4323
$ws.Range("N55").Value = -19354.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 16247.5
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 31495
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 31495
$ws.Range("M43").Value = -849
$ws.Range("N43").Value = -31797

$ws.Range("H46").Value = 34925
$ws.Range("J46").Value = 34925
$ws.Range("L46").Value = 34925
$ws.Range("N46").Value = -35237

$ws.Range("H57").Value = 14996.667
$ws.Range("J57").Value = 14996.667
$ws.Range("L57").Value = 14996.667
$ws.Range("N57").Value = -16636.667

$ws.Range("H80").Value = 7841.6665
$ws.Range("J80").Value = 10006
$ws.Range("L80").Value = 10006
$ws.Range("N80").Value = -12002

$ws.Range("H83").Value = 7841.6665
$ws.Range("J83").Value = 10006
$ws.Range("L83").Value = 50030
$ws.Range("N83").Value = -60014

$ws.Range("H97").Value = 2004.2
$ws.Range("I97").Value = 3070
$ws.Range("J97").Value = 405.5
$ws.Range("K97").Value = 3070
$ws.Range("L97").Value = 405.5
$ws.Range("M97").Value = -2574
$ws.Range("N97").Value = -1397.5

$ws.Range("H102").Value = 3590.4546
$ws.Range("I102").Value = 3849.5
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 3849.5
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -2227.5
$ws.Range("N102").Value = -4244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2629.76
$ws.Range("I132").Value = 1691.7059
$ws.Range("J132").Value = 4623.125
$ws.Range("K132").Value = 5075.1177
$ws.Range("L132").Value = 13869.375
$ws.Range("M132").Value = -2545.1177
$ws.Range("N132").Value = -18929.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5393.4546
$ws.Range("I81").Value = 1924.8889
$ws.Range("K81").Value = 3849.7778
$ws.Range("M81").Value = -2788.7778

$ws.Range("H84").Value = 5393.4546
$ws.Range("I84").Value = 1924.8889
$ws.Range("K84").Value = 19248.889
$ws.Range("M84").Value = -13944.889

$ws.Range("H107").Value = 623.9091
$ws.Range("J107").Value = 417.25
$ws.Range("L107").Value = 1251.75
$ws.Range("N107").Value = -5091.75

$ws.Range("H136").Value = 644.6667
$ws.Range("I136").Value = 634
$ws.Range("K136").Value = 1902
$ws.Range("M136").Value = 648
